$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the "test4" environment references with "test9" across the
# relevant cells (A2, C2, D2, F2, G2), preserving existing hyperlinks.
$ws.Range("A2").Value = "https://test9.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://test9.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://test9.cliotest.com/warehouse/control/main"
$ws.Range("F2").Value = "virtual_cabitest9"
$ws.Range("G2").Value = "test9"
